$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the testing_criterion default value (column M, row 2) from 8 to 16
$ws.Range("M2").Value = 16

# Update the active selection to N2 (matches the author's final cursor position)
$ws.Range("N2").Select()
